$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 27988.408
$ws.Range("I11").Value = 27988.408
$ws.Range("K11").Value = 27988.408
$ws.Range("M11").Value = -27848.408
$ws.Range("H28").Value = 665.6429000000001
$ws.Range("I28").Value = 601.46155
$ws.Range("K28").Value = 601.46155
$ws.Range("M28").Value = -116.46155
$ws.Range("H40").Value = 1900
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").Value = ""
$ws.Range("H51").Value = 3303.4634
$ws.Range("H80").Value = 1660.3704
$ws.Range("J80").Value = 1472.6842
$ws.Range("L80").Value = 4418.0526
$ws.Range("N80").Value = -6414.0526
$ws.Range("H83").Value = 1660.3704
$ws.Range("J83").Value = 1472.6842
$ws.Range("L83").Value = 13254.1578
$ws.Range("N83").Value = -23238.1578
$ws.Range("H87").Value = 62000
$ws.Range("J87").Value = 62000
$ws.Range("L87").Value = 62000
$ws.Range("N87").Value = -64496
$ws.Range("H90").Value = 62000
$ws.Range("J90").Value = 62000
$ws.Range("L90").Value = 186000
$ws.Range("N90").Value = -198480
$ws.Range("H94").Value = 1499
$ws.Range("I94").Value = 1499
$ws.Range("K94").Value = 1499
$ws.Range("M94").Value = -1048
$ws.Range("H99").Value = 322.33334
$ws.Range("I99").Value = 322.33334
$ws.Range("K99").Value = 967.0000200000001
$ws.Range("M99").Value = 530.9999799999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 74677.5
$ws.Range("J24").Value = 74677.5
$ws.Range("L24").Value = 74677.5
$ws.Range("N24").Value = -75425.5
$ws.Range("H61").Value = 4767093.5
$ws.Range("I61").Value = 5703.8237
$ws.Range("K61").Value = 5703.8237
$ws.Range("M61").Value = -5491.8237
$ws.Range("H88").Value = 1939.7084
$ws.Range("I88").Value = 1903.7142
$ws.Range("J88").Value = 1954.5294
$ws.Range("K88").Value = 1903.7142
$ws.Range("L88").Value = 1954.5294
$ws.Range("M88").Value = -1497.7142
$ws.Range("N88").Value = -2766.5294
$ws.Range("H91").Value = 1939.7084
$ws.Range("I91").Value = 1903.7142
$ws.Range("J91").Value = 1954.5294
$ws.Range("K91").Value = 1903.7142
$ws.Range("L91").Value = 1954.5294
$ws.Range("M91").Value = -499.7141999999999
$ws.Range("N91").Value = -4762.529399999999
$ws.Range("H100").Value = 74677.5
$ws.Range("J100").Value = 74677.5
$ws.Range("L100").Value = 74677.5
$ws.Range("N100").Value = -76841.5
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = ""
$ws.Range("H136").Value = 4767093.5
$ws.Range("I136").Value = 5703.8237
$ws.Range("K136").Value = 17111.4711
$ws.Range("M136").Value = -14561.4711

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3537
$ws.Range("I86").Value = 2847
$ws.Range("K86").Value = 2847
$ws.Range("M86").Value = -1724
$ws.Range("H89").Value = 3537
$ws.Range("I89").Value = 2847
$ws.Range("K89").Value = 14235
$ws.Range("M89").Value = -8619
$ws.Range("H107").Value = 1146
$ws.Range("I107").Value = 1146
$ws.Range("K107").Value = 1146
$ws.Range("M107").Value = 774
$ws.Range("H134").Value = 8774926
$ws.Range("I134").Value = 2759.7646
$ws.Range("K134").Value = 8279.293799999999
$ws.Range("M134").Value = -5744.293799999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2225147.2
$ws.Range("J31").Value = 2840.6155
$ws.Range("L31").Value = 2840.6155
$ws.Range("N31").Value = -3430.6155
$ws.Range("H34").Value = 2225147.2
$ws.Range("J34").Value = 2840.6155
$ws.Range("L34").Value = 2840.6155
$ws.Range("N34").Value = -3244.6155
$ws.Range("H62").Value = 4342
$ws.Range("I62").Value = 3995.5
$ws.Range("J62").Value = 4515.25
$ws.Range("K62").Value = 3995.5
$ws.Range("L62").Value = 4515.25
$ws.Range("M62").Value = -3371.5
$ws.Range("N62").Value = -5763.25
$ws.Range("H65").Value = 4342
$ws.Range("I65").Value = 3995.5
$ws.Range("J65").Value = 4515.25
$ws.Range("K65").Value = 19977.5
$ws.Range("L65").Value = 22576.25
$ws.Range("M65").Value = -16857.5
$ws.Range("N65").Value = -28816.25
$ws.Range("H132").Value = 2285.6956
$ws.Range("I132").Value = 2003.0526
$ws.Range("K132").Value = 6009.1578
$ws.Range("M132").Value = -3479.1578

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 140.86667
$ws.Range("I2").Value = 117.181816
$ws.Range("J2").Value = 206
$ws.Range("K2").Value = 703.0908959999999
$ws.Range("L2").Value = 1236
$ws.Range("M2").Value = -590.0908959999999
$ws.Range("N2").Value = -1462
$ws.Range("H17").Value = 593.75
$ws.Range("I17").Value = 182.66667
$ws.Range("K17").Value = 548.00001
$ws.Range("M17").Value = -379.00001
$ws.Range("H22").Value = 277.57144
$ws.Range("J22").Value = 501
$ws.Range("L22").Value = 1503
$ws.Range("N22").Value = -1841
$ws.Range("H27").Value = 277.57144
$ws.Range("J27").Value = 501
$ws.Range("L27").Value = 1503
$ws.Range("N27").Value = -1707
$ws.Range("H32").Value = 4024.75
$ws.Range("I32").Value = 3666.3333
$ws.Range("K32").Value = 10998.9999
$ws.Range("M32").Value = -10715.9999
$ws.Range("H46").Value = 1964.24
$ws.Range("I46").Value = 1671.9524
$ws.Range("K46").Value = 5015.857199999999
$ws.Range("M46").Value = -4924.857199999999
$ws.Range("H118").Value = 982
$ws.Range("I118").Value = 982
$ws.Range("K118").Value = 2946
$ws.Range("M118").Value = -1703
$ws.Range("H131").Value = 2030.4
$ws.Range("J131").Value = 2055.5
$ws.Range("L131").Value = 6166.5
$ws.Range("N131").Value = -16246.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2000
$ws.Range("J80").Value = 2000
$ws.Range("L80").Value = 2000
$ws.Range("N80").Value = -3996
$ws.Range("H83").Value = 2000
$ws.Range("J83").Value = 2000
$ws.Range("L83").Value = 10000
$ws.Range("N83").Value = -19984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1664.3793
$ws.Range("J55").Value = 1468.3158
$ws.Range("L55").Value = 1468.3158
$ws.Range("N55").Value = -1814.3158
$ws.Range("H82").Value = 1129.579
$ws.Range("I82").Value = 1217.6666
$ws.Range("K82").Value = 1217.6666
$ws.Range("M82").Value = -856.6666
$ws.Range("H85").Value = 1129.579
$ws.Range("I85").Value = 1217.6666
$ws.Range("K85").Value = 1217.6666
$ws.Range("M85").Value = 30.33339999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2107.818
$ws.Range("I81").Value = 1493.2632
$ws.Range("K81").Value = 2986.5264
$ws.Range("M81").Value = -1925.5264
$ws.Range("H84").Value = 2107.818
$ws.Range("I84").Value = 1493.2632
$ws.Range("K84").Value = 14932.632
$ws.Range("M84").Value = -9628.632000000001
$ws.Range("H132").Value = 11907549
$ws.Range("I132").Value = 15154005
$ws.Range("K132").Value = 45462015
$ws.Range("M132").Value = -45459485
$ws.Range("H136").Value = 18622688
$ws.Range("I136").Value = 10777896
$ws.Range("J136").Value = 50001860
$ws.Range("K136").Value = 32333688
$ws.Range("L136").Value = 150005580
$ws.Range("M136").Value = -32331138
